{"js": "// Edit 1: \"Metaclass, Class, Instance, Context, Occurrence, Role.\" -> \"...Role Resource Metada.\"\n// (the first, bullet-list occurrence near the top of the document)\n{\n  const results = context.document.body.search(\n    \"Metaclass, Class, Instance, Context, Occurrence, Role.\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"Metaclass, Class, Instance, Context, Occurrence, Role Resource Metada.\",\n      \"Replace\"\n    );\n    await context.sync();\n  }\n}\n\n// Edit 2: \"Sample Workflow: ToDo\" -> \"Sample Workflow:\" and insert 7 new bullet\n// paragraphs (same list formatting) right after it.\n{\n  const results = context.document.body.search(\"Sample Workflow: ToDo\", {\n    matchCase: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    const target = results.items[0];\n    const paragraph = target.paragraphs.getFirst();\n    await context.sync();\n\n    // Shrink the text of the found paragraph first.\n    target.insertText(\"Sample Workflow:\", \"Replace\");\n    await context.sync();\n\n    const newLines = [\n      \"aResource.flatMap(anStatement) : aResourceOccurrence;\",\n      \"aResourceOccurrence.flatMap(Activation::KindsCase) : aKindResource;\",\n      \"Activation::[Role]OccurrencesCase *: Occurrence[Role][];\",\n      \"DCI / MVC DDD Application Layer: OGM (Sesame Elmo / Alibaba. Qi4j). Core / Domains ontologies.\",\n      \"Resource: Types hierarchies / instances / occurrences (URN, Statement, CSPORole, Kind).\",\n      \"Relationships (discrete / continuous). Order. Translation / Equivalences entailments.\",\n      \"Metaclass, Class, Instance, Context, Occurrence, Role Resource Metadata Maps Monad with contextual CSPOs Statements (schema and occurrences) for Resources in Roles.\",\n    ];\n\n    let anchor = paragraph;\n    for (const line of newLines) {\n      anchor = anchor.insertParagraph(line, \"After\");\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Edit 1: \"Metaclass, Class, Instance, Context, Occurrence, Role.\" ->\n# \"...Role Resource Metada.\" (first, bullet-list occurrence near the top).\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Metaclass, Class, Instance, Context, Occurrence, Role.\") {\n        $p.Range.Text = \"Metaclass, Class, Instance, Context, Occurrence, Role Resource Metada.\"\n        break\n    }\n}\n\n# Edit 2: \"Sample Workflow: ToDo\" -> \"Sample Workflow:\" and insert 7 new bullet\n# paragraphs (same list formatting) right after it.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Sample Workflow: ToDo\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Text = \"Sample Workflow:\"\n\n    $newLines = @(\n        \"aResource.flatMap(anStatement) : aResourceOccurrence;\",\n        \"aResourceOccurrence.flatMap(Activation::KindsCase) : aKindResource;\",\n        \"Activation::[Role]OccurrencesCase *: Occurrence[Role][];\",\n        \"DCI / MVC DDD Application Layer: OGM (Sesame Elmo / Alibaba. Qi4j). Core / Domains ontologies.\",\n        \"Resource: Types hierarchies / instances / occurrences (URN, Statement, CSPORole, Kind).\",\n        \"Relationships (discrete / continuous). Order. Translation / Equivalences entailments.\",\n        \"Metaclass, Class, Instance, Context, Occurrence, Role Resource Metadata Maps Monad with contextual CSPOs Statements (schema and occurrences) for Resources in Roles.\"\n    )\n\n    $anchor = $target\n    foreach ($line in $newLines) {\n        $anchor.Range.InsertParagraphAfter()\n        $anchor = $anchor.Next()\n        $anchor.Range.Text = $line\n    }\n}\n\nWrite-Output \"done\"\n"}
